$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph that currently reads "4 black socks" and the
#    three sub-bullet paragraphs that directly follow it
#    ("Sub-problems" / "Selection done in complete darkness" /
#    "Socks cannot be seen until after selection is made").  These three
#    trailing paragraphs are being removed entirely, and the bookmark
#    that used to sit at the very end of the document (on the last of
#    those paragraphs) moves into the "socks" paragraph.
# ------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "4 black socks") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    throw "Could not find paragraph '4 black socks'"
}

# The three paragraphs to delete are the ones immediately after it.
$delStart = $d.Paragraphs.Item($targetIndex + 1)
$delEnd = $d.Paragraphs.Item($targetIndex + 3)

# Sanity-check their text before blowing them away.
$expectedTexts = @("Sub-problems", "Selection done in complete darkness", "Socks cannot be seen until after selection is made")
for ($j = 0; $j -lt 3; $j++) {
    $p = $d.Paragraphs.Item($targetIndex + 1 + $j)
    $txt = $p.Range.Text.TrimEnd()
    if ($txt -ne $expectedTexts[$j]) {
        throw "Unexpected paragraph text while deleting: '$txt'"
    }
}

$delRange = $d.Range($delStart.Range.Start, $delEnd.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 2. Change "black" to "white" inside the "4 black socks" paragraph.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$pStart = $p.Range.Start
$blackStart = $pStart + 2
$blackEnd = $blackStart + 5
$blackRange = $d.Range($blackStart, $blackEnd)
if ($blackRange.Text -ne "black") {
    throw "Expected 'black' at computed offset, found '$($blackRange.Text)'"
}
$blackRange.Text = "white"

# ------------------------------------------------------------------
# 3. Split "4 " away from "white" into separate runs. We do this by
#    wrapping "white" in a throw-away bookmark (which forces a run
#    split at both of its edges) and then deleting that bookmark again
#    -- the run split persists even after the bookmark itself is gone.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$pStart = $p.Range.Start
$whiteStart = $pStart + 2
$whiteEnd = $whiteStart + 5
$whiteRange = $d.Range($whiteStart, $whiteEnd)
if ($whiteRange.Text -ne "white") {
    throw "Expected 'white' at computed offset, found '$($whiteRange.Text)'"
}
$d.Bookmarks.Add("TempSplitMark", $whiteRange)
$d.Bookmarks.Item("TempSplitMark").Delete()

# ------------------------------------------------------------------
# 4. Re-insert the real "_GoBack" bookmark (collapsed / zero-length) at
#    the boundary between "white" and " socks", matching where it used
#    to sit relative to the final edit.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$pStart = $p.Range.Start
$afterWhite = $pStart + 7
$bmRange = $d.Range($afterWhite, $afterWhite)
$d.Bookmarks.Add("_GoBack", $bmRange)

$finalText = $d.Paragraphs.Item($targetIndex).Range.Text
Write-Output "Final paragraph text: $finalText"
